$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 45 (La Araucania / Poroto
# granado, "Region del Maule" origin, 2023-12-12). This pushes every
# existing record that used to live in rows 45-123 down by one row, to
# rows 46-124, without altering their values.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(45,1).Value = 10
$ws.Cells.Item(45,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45,3).Value = "La Araucanía"
$ws.Cells.Item(45,4).Value = "2023-12-12"
$ws.Cells.Item(45,5).Value = 9
$ws.Cells.Item(45,6).Value = 100112030
$ws.Cells.Item(45,7).Value = "Poroto granado"
$ws.Cells.Item(45,8).Value = "Sin especificar"
$ws.Cells.Item(45,9).Value = "Primera"
$ws.Cells.Item(45,10).Value = 15
$ws.Cells.Item(45,11).Value = 60000
$ws.Cells.Item(45,12).Value = 60000
$ws.Cells.Item(45,13).Value = 60000
$ws.Cells.Item(45,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45,15).Value = "Región del Maule"
$ws.Cells.Item(45,16).Value = 2400
$ws.Cells.Item(45,17).Value = 25
$ws.Cells.Item(45,18).Value = "Hortaliza"
